$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells D1:I1, copying style of existing header row (A1:C1)
$ws.Range("A1:C1").Copy() | Out-Null
$ws.Range("D1:I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("D1").Value = "Stat - alpha"
$ws.Range("E1").Value = "Stat - epsilon"
$ws.Range("F1").Value = "Stat - episodes"
$ws.Range("G1").Value = "Hyper - alpha"
$ws.Range("H1").Value = "Hyper - epsilon"
$ws.Range("I1").Value = "Hyper - episodes"

# Update existing rows' timestamps / stats
$ws.Range("A2").Value = "2025-06-15 19:08:20"
$ws.Range("C2").Value = "{'win': 1, 'loss': 0, 'draw': 1}"

$ws.Range("A3").Value = "2025-06-15 19:08:24"
$ws.Range("C3").Value = "{'win': 2, 'loss': 0, 'draw': 2}"

$ws.Range("A4").Value = "2025-06-15 19:08:27"
$ws.Range("C4").Value = "{'win': 2, 'loss': 0, 'draw': 4}"

$ws.Range("A5").Value = "2025-06-15 19:57:47"
$ws.Range("C5").Value = "{'win': 2, 'loss': 0, 'draw': 0}"
$ws.Range("D5").Value = 0.1
$ws.Range("E5").Value = 0.1
$ws.Range("F5").Value = 1000
$ws.Range("G5").Value = 0.1
$ws.Range("H5").Value = 0.1
$ws.Range("I5").Value = 1000

$ws.Range("A6").Value = "2025-06-15 19:57:54"
$ws.Range("C6").Value = "{'win': 2, 'loss': 0, 'draw': 2}"
$ws.Range("D6").Value = 0.1
$ws.Range("E6").Value = 0.1
$ws.Range("F6").Value = 1000
$ws.Range("G6").Value = 0.1
$ws.Range("H6").Value = 0.1
$ws.Range("I6").Value = 1000

$ws.Range("A7").Value = "2025-06-15 19:57:59"
$ws.Range("C7").Value = "{'win': 2, 'loss': 0, 'draw': 4}"
$ws.Range("D7").Value = 0.1
$ws.Range("E7").Value = 0.1
$ws.Range("F7").Value = 1000
$ws.Range("G7").Value = 0.1
$ws.Range("H7").Value = 0.1
$ws.Range("I7").Value = 1000

$ws.Range("A8").Value = "2025-06-15 19:58:03"
$ws.Range("C8").Value = "{'win': 3, 'loss': 0, 'draw': 5}"
$ws.Range("D8").Value = 0.1
$ws.Range("E8").Value = 0.1
$ws.Range("F8").Value = 1000
$ws.Range("G8").Value = 0.1
$ws.Range("H8").Value = 0.1
$ws.Range("I8").Value = 1000

# Remove row 9 entirely
$ws.Rows.Item(9).Delete() | Out-Null

# Touch D2:I4 so an (empty) cell record is emitted, matching upstream export
$ws.Range("D2:I4").WrapText = $false
